$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 4-7 were blank filler rows; bring in the same formatting (styles + row
# height) already used by the data row above (row 3), then fill in the new
# judge non-availability records.
$ws.Range("A3:I3").Copy()
$ws.Range("A4:I7").PasteSpecial(-4122)
$ws.Rows.Item(4).RowHeight = 35
$ws.Rows.Item(5).RowHeight = 35
$ws.Rows.Item(6).RowHeight = 35
$ws.Rows.Item(7).RowHeight = 35

# Row 4: Greneven, Eliana - Central
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = 43205
$ws.Range("C4").Value = "Central"
$ws.Range("E4").Value = "2 (1W200B)"
$ws.Range("G4").Value = "Greneven, Eliana "
$ws.Range("F4").Value = "BVAGRENEV"
$ws.Range("H4").Value = 1

# Row 5: Mulligan, James - Central
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = 43205
$ws.Range("C5").Value = "Central"
$ws.Range("E5").Value = "2 (1W200B)"
$ws.Range("F5").Value = "BVAMULLIGAN"
$ws.Range("G5").Value = "Mulligan, James"
$ws.Range("H5").Value = 1

# Row 6: Mulligan, James - Video
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = 43205
$ws.Range("C6").Value = "Video"
$ws.Range("E6").Value = "2 (1W200B)"
$ws.Range("F6").Value = "BVAMULLIGAN"
$ws.Range("G6").Value = "Mulligan, James"
$ws.Range("H6").Value = 1

# Row 7: Mulligan, James - Virtual
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = 43205
$ws.Range("C7").Value = "Virtual"
$ws.Range("E7").Value = "2 (1W200B)"
$ws.Range("F7").Value = "BVAMULLIGAN"
$ws.Range("G7").Value = "Mulligan, James"
$ws.Range("H7").Value = 1

$ws.Range("H18").Select() | Out-Null
